$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 6: Days of Chunder | Antidote
$ws.Range("H6").Value = 2000050
$ws.Range("J6").Value = 49
$ws.Range("L6").Value = 147
$ws.Range("N6").Value = -371

# Row 62: The Mustache Suits Him | Enchanted Mythrite Ink
$ws.Range("H62").Value = 5629.923
$ws.Range("I62").Value = 5919
$ws.Range("J62").Value = 4666.3335
$ws.Range("K62").Value = 5919
$ws.Range("L62").Value = 4666.3335
$ws.Range("M62").Value = -5295
$ws.Range("N62").Value = -5914.3335

# Row 65: Forgery of Convenience (L) | Enchanted Mythrite Ink
$ws.Range("H65").Value = 5629.923
$ws.Range("I65").Value = 5919
$ws.Range("J65").Value = 4666.3335
$ws.Range("K65").Value = 29595
$ws.Range("L65").Value = 23331.6675
$ws.Range("M65").Value = -26475
$ws.Range("N65").Value = -29571.6675

# Row 88: The Grave of Hemlock Groves | Growth Formula Zeta
$ws.Range("H88").Value = 1267.3636
$ws.Range("I88").Value = 998
$ws.Range("J88").Value = 1294.3
$ws.Range("K88").Value = 998
$ws.Range("L88").Value = 1294.3
$ws.Range("M88").Value = -592
$ws.Range("N88").Value = -2106.3

# Row 91: Dappling the Highlands (L) | Growth Formula Zeta
$ws.Range("H91").Value = 1267.3636
$ws.Range("I91").Value = 998
$ws.Range("J91").Value = 1294.3
$ws.Range("K91").Value = 998
$ws.Range("L91").Value = 1294.3
$ws.Range("M91").Value = 406
$ws.Range("N91").Value = -4102.3

# Row 98: The Dotted Line | Enchanted Durium Ink
$ws.Range("H98").Value = 795.2857
$ws.Range("I98").Value = 778.92
$ws.Range("K98").Value = 778.92
$ws.Range("M98").Value = 719.08

# Row 106: Making Your Mark | Enchanted Palladium Ink
$ws.Range("H106").Value = 1699.5
$ws.Range("I106").Value = 1699.5
$ws.Range("K106").Value = 1699.5
$ws.Range("M106").Value = -1068.5

# Row 112: Making Ends Meet | Superior Spiritbond Potion
$ws.Range("H112").Value = 60749.41
$ws.Range("I112").Value = 1825
$ws.Range("K112").Value = 5475
$ws.Range("M112").Value = -4367

# Row 122: Wishful Inking | Enchanted High Durium Ink
$ws.Range("H122").Value = 795.2857
$ws.Range("I122").Value = 778.92
$ws.Range("K122").Value = 2336.76
$ws.Range("M122").Value = 113.2400000000002

# Row 132: Fast-forwarding Flora | Growth Formula Lambda
$ws.Range("H132").Value = 2017.5667
$ws.Range("I132").Value = 1423
$ws.Range("K132").Value = 4269
$ws.Range("M132").Value = -1739

# Row 137: Cutting Edge of Culinary Quality | Magnesia Whetstone
$ws.Range("H137").Value = 6087.048
$ws.Range("I137").Value = 1271.7142
$ws.Range("J137").Value = 15717.714
$ws.Range("K137").Value = 3815.1426
$ws.Range("L137").Value = 47153.142
$ws.Range("M137").Value = -1265.1426
$ws.Range("N137").Value = -52253.142

# Row 138: All-night Crafting | Cunning Craftsman's Tisane
$ws.Range("H138").Value = 26875.637
$ws.Range("I138").Value = 57903.11
$ws.Range("J138").Value = 5395.077
$ws.Range("K138").Value = 173709.33
$ws.Range("L138").Value = 16185.231
$ws.Range("M138").Value = -168569.33
$ws.Range("N138").Value = -26465.231

$ws = $wb.Worksheets.Item("ARM")
# Row 32: Ingot We Trust | Steel Ingot
$ws.Range("H32").Value = 5187.5713
$ws.Range("I32").Value = 4694.237
$ws.Range("K32").Value = 4694.237
$ws.Range("M32").Value = -4407.237

# Row 132: Don't Bore Me, Ore Me | Mountain Chromite Ingot
$ws.Range("H132").Value = 3214.4722
$ws.Range("I132").Value = 3135.1738
$ws.Range("J132").Value = 3354.7693
$ws.Range("K132").Value = 9405.5214
$ws.Range("L132").Value = 10064.3079
$ws.Range("M132").Value = -6875.5214
$ws.Range("N132").Value = -15124.3079

# Row 135: Forgiveness for My Shins | Ruthenium Sabatons of Fending
$ws.Range("H135").Value = 141809.33
$ws.Range("J135").Value = 141809.33
$ws.Range("L135").Value = 141809.33
$ws.Range("N135").Value = -151949.33

$ws = $wb.Worksheets.Item("BSM")
# Row 11: Down on the Pharm | Amateur's Mortar
$ws.Range("H11").Value = 11927.5
$ws.Range("I11").Value = 23420
$ws.Range("J11").Value = 435
$ws.Range("K11").Value = 23420
$ws.Range("L11").Value = 435
$ws.Range("M11").Value = -23280
$ws.Range("N11").Value = -715

# Row 96: Hammer Time | High Steel Sledgehammer
$ws.Range("H96").Value = 16662.334
$ws.Range("I96").Value = 13106
$ws.Range("K96").Value = 13106
$ws.Range("M96").Value = -10360

# Row 105: Ingot to Wing It | Molybdenum Ingot
$ws.Range("H105").Value = 2057.4
$ws.Range("I105").Value = 1508.7646
$ws.Range("K105").Value = 1508.7646
$ws.Range("M105").Value = 238.2354

# Row 134: Ruthenium Supremium | Ruthenium Ingot
$ws.Range("H134").Value = 1818.1794
$ws.Range("I134").Value = 1813.3948
$ws.Range("K134").Value = 5440.1844
$ws.Range("M134").Value = -2905.1844

$ws = $wb.Worksheets.Item("CRP")
# Row 3: Touch and Heal | Maple Pattens
$ws.Range("H3").Value = 4150
$ws.Range("I3").Value = 5575
$ws.Range("K3").Value = 5575
$ws.Range("M3").Value = -5462

# Row 8: Bows for the Boys | Maple Longbow
$ws.Range("H8").Value = 500
$ws.Range("I8").Value = 787.8
$ws.Range("J8").Value = 212.2
$ws.Range("K8").Value = 787.8
$ws.Range("L8").Value = 212.2
$ws.Range("M8").Value = -647.8
$ws.Range("N8").Value = -492.2

# Row 31: Wall Not Found | Walnut Lumber
$ws.Range("H31").Value = 26486.643
$ws.Range("I31").Value = 40507.19
$ws.Range("K31").Value = 40507.19
$ws.Range("M31").Value = -40212.19

# Row 34: Armoires of the Rich and Famous | Walnut Lumber
$ws.Range("H34").Value = 26486.643
$ws.Range("I34").Value = 40507.19
$ws.Range("K34").Value = 40507.19
$ws.Range("M34").Value = -40305.19

# Row 53: A Winning Combo | Oak Composite Bow
$ws.Range("H53").Value = 24999.5
$ws.Range("J53").Value = 24999.5
$ws.Range("L53").Value = 24999.5
$ws.Range("N53").Value = -26213.5

# Row 58: You Do the Heavy Lifting | Mahogany Lumber
$ws.Range("H58").Value = 2762.65
$ws.Range("J58").Value = 3022.5833
$ws.Range("L58").Value = 3022.5833
$ws.Range("N58").Value = -3428.5833

# Row 74: License to Heal | Dark Chestnut Rod
$ws.Range("H74").Value = 22500
$ws.Range("J74").Value = 22500
$ws.Range("L74").Value = 22500
$ws.Range("N74").Value = -24248

# Row 77: Purified Polyrhythm (L) | Dark Chestnut Rod
$ws.Range("H77").Value = 22500
$ws.Range("J77").Value = 22500
$ws.Range("L77").Value = 67500
$ws.Range("N77").Value = -76236

# Row 99: O Pine | Pine Lumber
$ws.Range("H99").Value = 7629.75
$ws.Range("I99").Value = 2752.5
$ws.Range("J99").Value = 12507
$ws.Range("K99").Value = 2752.5
$ws.Range("L99").Value = 12507
$ws.Range("M99").Value = -1254.5
$ws.Range("N99").Value = -15503

# Row 122: Timber of Tenkonto | Horse Chestnut Lumber
$ws.Range("H122").Value = 1919.125
$ws.Range("I122").Value = 1809
$ws.Range("J122").Value = 2249.5
$ws.Range("K122").Value = 5427
$ws.Range("L122").Value = 6748.5
$ws.Range("M122").Value = -2977
$ws.Range("N122").Value = -11648.5

# Row 126: A Better Conductor | Red Pine Lumber
$ws.Range("H126").Value = 7629.75
$ws.Range("I126").Value = 2752.5
$ws.Range("J126").Value = 12507
$ws.Range("K126").Value = 8257.5
$ws.Range("L126").Value = 37521
$ws.Range("M126").Value = -5787.5
$ws.Range("N126").Value = -42461

# Row 134: Wood You Be Quiet | Ceiba Lumber
$ws.Range("H134").Value = 23895.256
$ws.Range("I134").Value = 18759.428
$ws.Range("K134").Value = 56278.284
$ws.Range("M134").Value = -53743.284

# Row 136: Turali Quality | Dark Mahogany Lumber
$ws.Range("H136").Value = 2762.65
$ws.Range("J136").Value = 3022.5833
$ws.Range("L136").Value = 9067.749899999999
$ws.Range("N136").Value = -14167.7499

$ws = $wb.Worksheets.Item("CUL")
# Row 26: A Grape Idea | Grape Juice
$ws.Range("H26").Value = 1028.8334
$ws.Range("I26").Value = 36.22222
$ws.Range("J26").Value = 4006.6667
$ws.Range("K26").Value = 108.66666
$ws.Range("L26").Value = 12020.0001
$ws.Range("M26").Value = 179.33334
$ws.Range("N26").Value = -12596.0001

# Row 122: Salt of the North | Northern Sea Salt
$ws.Range("H122").Value = 1071.12
$ws.Range("I122").Value = 1308.4615
$ws.Range("J122").Value = 814
$ws.Range("K122").Value = 11776.1535
$ws.Range("L122").Value = 7326
$ws.Range("M122").Value = -9326.153499999999
$ws.Range("N122").Value = -12226

# Row 131: The Mountain Steeped | Tsai tou Vounou
$ws.Range("H131").Value = 34502.355
$ws.Range("I131").Value = 100822.6
$ws.Range("K131").Value = 302467.8
$ws.Range("M131").Value = -297427.8

$ws = $wb.Worksheets.Item("GSM")
# Row 130: Planisphere to Paper | Chondrite Magitek Planisphere
$ws.Range("H130").Value = 64999
$ws.Range("J130").Value = 64999
$ws.Range("L130").Value = 64999
$ws.Range("N130").Value = -75039

# Row 132: On Board for Lar | Lar Ingot
$ws.Range("H132").Value = 2831.1316
$ws.Range("I132").Value = 2730.9429
$ws.Range("J132").Value = 4000
$ws.Range("K132").Value = 8192.8287
$ws.Range("L132").Value = 12000
$ws.Range("M132").Value = -5662.8287
$ws.Range("N132").Value = -17060

$ws = $wb.Worksheets.Item("LTW")
# Row 2: Red in the Head | Leather Calot
$ws.Range("H2").Value = 1916721.8
$ws.Range("I2").Value = 350166.34
$ws.Range("J2").Value = 2699999.5
$ws.Range("K2").Value = 350166.34
$ws.Range("L2").Value = 2699999.5
$ws.Range("M2").Value = -350054.34
$ws.Range("N2").Value = -2700223.5

# Row 46: Supply Side Logic | Boar Leather
$ws.Range("H46").Value = 1031.4445
$ws.Range("I46").Value = 1047.875
$ws.Range("K46").Value = 1047.875
$ws.Range("M46").Value = -859.875

# Row 68: You Could Say It's a Moving Target | Wyvern Leather
$ws.Range("H68").Value = 2922.5
$ws.Range("I68").Value = 2922.5
$ws.Range("K68").Value = 2922.5
$ws.Range("M68").Value = -2173.5

# Row 71: They Call It Bloody Mary (L) | Wyvern Leather
$ws.Range("H71").Value = 2922.5
$ws.Range("I71").Value = 2922.5
$ws.Range("K71").Value = 14612.5
$ws.Range("M71").Value = -10868.5

# Row 132: Tenets of Tanning | Silver Lobo Leather
$ws.Range("H132").Value = 4096.773
$ws.Range("I132").Value = 3979.111
$ws.Range("K132").Value = 11937.333
$ws.Range("M132").Value = -9407.332999999999

# Row 136: Respect for Br'aax | Br'aax Leather
$ws.Range("H136").Value = 35365.1
$ws.Range("I136").Value = 40105.883
$ws.Range("J136").Value = 4550
$ws.Range("K136").Value = 120317.649
$ws.Range("L136").Value = 13650
$ws.Range("M136").Value = -117767.649
$ws.Range("N136").Value = -18750

$ws = $wb.Worksheets.Item("WVR")
# Row 3: Trew Enough | Hempen Chausses
$ws.Range("H3").Value = 85643.586
$ws.Range("I3").Value = 126016.875
$ws.Range("J3").Value = 4897
$ws.Range("K3").Value = 126016.875
$ws.Range("L3").Value = 4897
$ws.Range("M3").Value = -125902.875
$ws.Range("N3").Value = -5125

# Row 11: Wiggle Room | Hempen Shepherd's Tunic
$ws.Range("H11").Value = 170164
$ws.Range("J11").Value = 4196.8
$ws.Range("L11").Value = 4196.8
$ws.Range("N11").Value = -4480.8

# Row 113: A Tender Table | Pixie Floss
$ws.Range("H113").Value = 1026.0714
$ws.Range("I113").Value = 791.625
$ws.Range("J113").Value = 1338.6666
$ws.Range("K113").Value = 2374.875
$ws.Range("L113").Value = 4015.9998
$ws.Range("M113").Value = -204.875
$ws.Range("N113").Value = -8355.9998

# Row 118: Something in My Eye | Ovim Wool Turban of Gathering
$ws.Range("H118").Value = 0
$ws.Range("J118").Value = 0
$ws.Range("L118").Value = 0
$ws.Range("N118").ClearContents()

# Row 132: Comfy Cabins | Snow Cotton Cloth
$ws.Range("H132").Value = 1035.0834
$ws.Range("I132").Value = 993.13043
$ws.Range("K132").Value = 2979.39129
$ws.Range("M132").Value = -449.39129

# Row 135: In Line with Linen | Mountain Linen Cloak of Casting
$ws.Range("H135").Value = 35000
$ws.Range("J135").Value = 35000
$ws.Range("L135").Value = 35000
$ws.Range("N135").Value = -45140
